$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 243.55556
$ws.Range("I4").Value = 224
$ws.Range("K4").Value = 224
$ws.Range("M4").Value = -110
$ws.Range("H12").Value = 181.75
$ws.Range("I12").Value = 157.8
$ws.Range("J12").Value = 301.5
$ws.Range("K12").Value = 157.8
$ws.Range("L12").Value = 301.5
$ws.Range("M12").Value = 12.19999999999999
$ws.Range("N12").Value = -641.5
$ws.Range("H32").Value = 4391.25
$ws.Range("I32").Value = 4150
$ws.Range("J32").Value = 4471.6665
$ws.Range("K32").Value = 4150
$ws.Range("L32").Value = 4471.6665
$ws.Range("M32").Value = -3824
$ws.Range("N32").Value = -5123.6665
$ws.Range("H38").Value = 33.57143
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("H76").Value = 3281
$ws.Range("I76").Value = 3281
$ws.Range("K76").Value = 3281
$ws.Range("M76").Value = -2966
$ws.Range("H79").Value = 3281
$ws.Range("I79").Value = 3281
$ws.Range("K79").Value = 3281
$ws.Range("M79").Value = -2189
$ws.Range("H86").Value = 2567.7144
$ws.Range("J86").Value = 2691.7273
$ws.Range("L86").Value = 2691.7273
$ws.Range("N86").Value = -4937.7273
$ws.Range("H88").Value = 14891
$ws.Range("I88").Value = 3750
$ws.Range("K88").Value = 3750
$ws.Range("M88").Value = -3344
$ws.Range("H89").Value = 2567.7144
$ws.Range("J89").Value = 2691.7273
$ws.Range("L89").Value = 13458.6365
$ws.Range("N89").Value = -24690.6365
$ws.Range("H91").Value = 14891
$ws.Range("I91").Value = 3750
$ws.Range("K91").Value = 3750
$ws.Range("M91").Value = -2346
$ws.Range("H137").Value = 383802.34
$ws.Range("I137").Value = 2159.2258
$ws.Range("J137").Value = 602893.75
$ws.Range("K137").Value = 6477.6774
$ws.Range("L137").Value = 1808681.25
$ws.Range("M137").Value = -3927.6774
$ws.Range("N137").Value = -1813781.25
$ws.Range("H141").Value = 1887.5454
$ws.Range("I141").Value = 1711
$ws.Range("K141").Value = 5133
$ws.Range("M141").Value = 47
$ws.Range("N38").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9808632
$ws.Range("I32").Value = 10421646
$ws.Range("J32").Value = 394.5
$ws.Range("K32").Value = 10421646
$ws.Range("L32").Value = 394.5
$ws.Range("M32").Value = -10421359
$ws.Range("H45").Value = 2239.5
$ws.Range("I45").Value = 2251
$ws.Range("K45").Value = 2251
$ws.Range("M45").Value = -1874
$ws.Range("H61").Value = 5573305.5
$ws.Range("I61").Value = 5573305.5
$ws.Range("K61").Value = 5573305.5
$ws.Range("M61").Value = -5573093.5
$ws.Range("H74").Value = 11796.333
$ws.Range("I74").Value = 13494.5
$ws.Range("J74").Value = 8400
$ws.Range("K74").Value = 13494.5
$ws.Range("L74").Value = 8400
$ws.Range("M74").Value = -12620.5
$ws.Range("N74").Value = -10148
$ws.Range("H77").Value = 11796.333
$ws.Range("I77").Value = 13494.5
$ws.Range("J77").Value = 8400
$ws.Range("K77").Value = 67472.5
$ws.Range("L77").Value = 42000
$ws.Range("M77").Value = -63104.5
$ws.Range("N77").Value = -50736
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H122").Value = 2472.5
$ws.Range("I122").Value = 2573.5881
$ws.Range("J122").Value = 1899.6666
$ws.Range("K122").Value = 7720.7643
$ws.Range("L122").Value = 5698.9998
$ws.Range("M122").Value = -5270.7643
$ws.Range("N122").Value = -10598.9998
$ws.Range("H132").Value = 989516.5600000001
$ws.Range("I132").Value = 1187199.6
$ws.Range("J132").Value = 1101.5
$ws.Range("K132").Value = 3561598.8
$ws.Range("L132").Value = 3304.5
$ws.Range("M132").Value = -3559068.8
$ws.Range("N132").Value = -8364.5
$ws.Range("H136").Value = 5573305.5
$ws.Range("I136").Value = 5573305.5
$ws.Range("K136").Value = 16719916.5
$ws.Range("M136").Value = -16717366.5
$ws.Range("N32").Value = -968.5
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 51783.25
$ws.Range("J2").Value = 51783.25
$ws.Range("L2").Value = 51783.25
$ws.Range("N2").Value = -52009.25
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("H20").Value = 1658.2142
$ws.Range("I20").Value = 1045.25
$ws.Range("J20").Value = 1903.4
$ws.Range("K20").Value = 1045.25
$ws.Range("L20").Value = 1903.4
$ws.Range("M20").Value = -798.25
$ws.Range("N20").Value = -2397.4
$ws.Range("H86").Value = 2343.7144
$ws.Range("I86").Value = 2401
$ws.Range("K86").Value = 2401
$ws.Range("M86").Value = -1278
$ws.Range("H89").Value = 2343.7144
$ws.Range("I89").Value = 2401
$ws.Range("K89").Value = 12005
$ws.Range("M89").Value = -6389
$ws.Range("H107").Value = 7341.9165
$ws.Range("I107").Value = 5935
$ws.Range("J107").Value = 8748.833000000001
$ws.Range("K107").Value = 5935
$ws.Range("L107").Value = 8748.833000000001
$ws.Range("M107").Value = -4015
$ws.Range("N107").Value = -12588.833
$ws.Range("H134").Value = 2187515.2
$ws.Range("I134").Value = 2980218
$ws.Range("K134").Value = 8940654
$ws.Range("M134").Value = -8938119
$ws.Range("N13").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 233279.06
$ws.Range("I31").Value = 348271.12
$ws.Range("J31").Value = 60791
$ws.Range("K31").Value = 348271.12
$ws.Range("L31").Value = 60791
$ws.Range("M31").Value = -347976.12
$ws.Range("N31").Value = -61381
$ws.Range("H34").Value = 233279.06
$ws.Range("I34").Value = 348271.12
$ws.Range("J34").Value = 60791
$ws.Range("K34").Value = 348271.12
$ws.Range("L34").Value = 60791
$ws.Range("M34").Value = -348069.12
$ws.Range("N34").Value = -61195
$ws.Range("H86").Value = 87739.5
$ws.Range("I86").Value = 4325.4614
$ws.Range("K86").Value = 4325.4614
$ws.Range("M86").Value = -3202.4614
$ws.Range("H89").Value = 87739.5
$ws.Range("I89").Value = 4325.4614
$ws.Range("K89").Value = 21627.307
$ws.Range("M89").Value = -16011.307
$ws.Range("H132").Value = 62827260
$ws.Range("I132").Value = 90913000
$ws.Range("K132").Value = 272739000
$ws.Range("M132").Value = -272736470

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 10770
$ws.Range("J82").Value = 17353.75
$ws.Range("L82").Value = 52061.25
$ws.Range("N82").Value = -52873.25
$ws.Range("H85").Value = 10770
$ws.Range("J85").Value = 17353.75
$ws.Range("L85").Value = 52061.25
$ws.Range("N85").Value = -54869.25
$ws.Range("H140").Value = 2564.7083
$ws.Range("I140").Value = 1847.0625
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 5541.1875
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = -361.1875
$ws.Range("N140").Value = -22360

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1200.1852
$ws.Range("I97").Value = 916.63635
$ws.Range("J97").Value = 2447.8
$ws.Range("K97").Value = 916.63635
$ws.Range("L97").Value = 2447.8
$ws.Range("M97").Value = -420.63635
$ws.Range("N97").Value = -3439.8
$ws.Range("H126").Value = 880271.5600000001
$ws.Range("I126").Value = 1854306.8
$ws.Range("J126").Value = 3639.9
$ws.Range("K126").Value = 5562920.4
$ws.Range("L126").Value = 10919.7
$ws.Range("M126").Value = -5560450.4
$ws.Range("N126").Value = -15859.7
$ws.Range("H132").Value = 34908580
$ws.Range("I132").Value = 48199990
$ws.Range("K132").Value = 144599970
$ws.Range("M132").Value = -144597440

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2349.3333
$ws.Range("I16").Value = 2149
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 2149
$ws.Range("L16").Value = 2750
$ws.Range("M16").Value = -1979
$ws.Range("N16").Value = -3090
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5272
$ws.Range("H136").Value = 47644.68
$ws.Range("I136").Value = 2695.625
$ws.Range("J136").Value = 107576.75
$ws.Range("K136").Value = 8086.875
$ws.Range("L136").Value = 322730.25
$ws.Range("M136").Value = -5536.875
$ws.Range("N136").Value = -327830.25
$ws.Range("M40").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1605.1765
$ws.Range("I107").Value = 1003.1923
$ws.Range("K107").Value = 3009.5769
$ws.Range("M107").Value = -1089.5769
$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 21000
$ws.Range("H132").Value = 3598033.5
$ws.Range("J132").Value = 7985.8887
$ws.Range("L132").Value = 23957.6661
$ws.Range("N132").Value = -29017.6661
$ws.Range("N126").Value = -25940
$ws.Range("M126").ClearContents()
